$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the tool names in column B (rows 2-10) with tool1/tool2/tool3 repeating
$ws.Range("B2").Value = "tool1"
$ws.Range("B3").Value = "tool2"
$ws.Range("B4").Value = "tool3"
$ws.Range("B5").Value = "tool1"
$ws.Range("B6").Value = "tool2"
$ws.Range("B7").Value = "tool3"
$ws.Range("B8").Value = "tool1"
$ws.Range("B9").Value = "tool2"
$ws.Range("B10").Value = "tool3"

# Update the active cell selection to B6
$ws.Range("B6").Select()
